$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 242, shifting existing rows 242:327 down to 243:328
$ws.Rows.Item(242).Insert()

# Populate the newly inserted row 242 with its data.
$ws.Range("A242").Value = 8
$ws.Range("B242").Value = "Terminal La Palmera de La Serena"
$ws.Range("C242").Value = "Coquimbo"
$ws.Range("D242").Value = 44900
$ws.Range("E242").Value = 4
$ws.Range("F242").Value = 100112021
$ws.Range("G242").Value = "Ají"
$ws.Range("H242").Value = "Inferno"
$ws.Range("I242").Value = "Primera"
$ws.Range("J242").Value = 460
$ws.Range("K242").Value = 12000
$ws.Range("L242").Value = 13000
$ws.Range("M242").Value = 12500
$ws.Range("N242").Value = "`$/caja 10 kilos"
$ws.Range("O242").Value = "Región de Arica y Parinacota"
$ws.Range("P242").Value = 1250
$ws.Range("Q242").Value = 10
$ws.Range("R242").Value = "Hortaliza"
